$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 761.3333
$ws.Cells.Item(38, 9).Value = 761.3333
$ws.Cells.Item(38, 11).Value = 2283.9999
$ws.Cells.Item(38, 13).Value = -1911.9999
$ws.Cells.Item(40, 8).Value = 2253.6155
$ws.Cells.Item(40, 9).Value = 1964.1428
$ws.Cells.Item(40, 10).Value = 2591.3333
$ws.Cells.Item(40, 11).Value = 1964.1428
$ws.Cells.Item(40, 12).Value = 2591.3333
$ws.Cells.Item(40, 13).Value = -1789.1428
$ws.Cells.Item(40, 14).Value = -2941.3333
$ws.Cells.Item(129, 8).Value = 3554.182
$ws.Cells.Item(129, 10).Value = 4000
$ws.Cells.Item(129, 12).Value = 12000
$ws.Cells.Item(129, 14).Value = -22000
$ws.Cells.Item(132, 8).Value = 1482.16
$ws.Cells.Item(132, 9).Value = 1472.45
$ws.Cells.Item(132, 10).Value = 1521
$ws.Cells.Item(132, 11).Value = 4417.35
$ws.Cells.Item(132, 12).Value = 4563
$ws.Cells.Item(132, 13).Value = -1887.35
$ws.Cells.Item(132, 14).Value = -9623
$ws.Cells.Item(135, 8).Value = 1477.7142
$ws.Cells.Item(135, 9).Value = 1104.8889
$ws.Cells.Item(135, 10).Value = 2148.8
$ws.Cells.Item(135, 11).Value = 9944.000099999999
$ws.Cells.Item(135, 12).Value = 19339.2
$ws.Cells.Item(135, 13).Value = -7409.000099999999
$ws.Cells.Item(135, 14).Value = -24409.2
$ws.Cells.Item(138, 8).Value = 5644.913
$ws.Cells.Item(138, 9).Value = 3473.75
$ws.Cells.Item(138, 10).Value = 5851.6904
$ws.Cells.Item(138, 11).Value = 10421.25
$ws.Cells.Item(138, 12).Value = 17555.0712
$ws.Cells.Item(138, 13).Value = -5281.25
$ws.Cells.Item(138, 14).Value = -27835.0712
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 19898.36
$ws.Cells.Item(32, 9).Value = 12117.944
$ws.Cells.Item(32, 10).Value = 24274.844
$ws.Cells.Item(32, 11).Value = 12117.944
$ws.Cells.Item(32, 12).Value = 24274.844
$ws.Cells.Item(32, 13).Value = -11830.944
$ws.Cells.Item(32, 14).Value = -24848.844
$ws.Cells.Item(45, 8).Value = 2031.1111
$ws.Cells.Item(45, 9).Value = 2010
$ws.Cells.Item(45, 11).Value = 2010
$ws.Cells.Item(45, 13).Value = -1633
$ws.Cells.Item(122, 8).Value = 591338.9
$ws.Cells.Item(122, 9).Value = 1113195.6
$ws.Cells.Item(122, 11).Value = 3339586.8
$ws.Cells.Item(122, 13).Value = -3337136.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2462.6
$ws.Cells.Item(86, 9).Value = 1376.5
$ws.Cells.Item(86, 11).Value = 1376.5
$ws.Cells.Item(86, 13).Value = -253.5
$ws.Cells.Item(89, 8).Value = 2462.6
$ws.Cells.Item(89, 9).Value = 1376.5
$ws.Cells.Item(89, 11).Value = 6882.5
$ws.Cells.Item(89, 13).Value = -1266.5
$ws.Cells.Item(99, 8).Value = 1000
$ws.Cells.Item(99, 9).Value = 1000
$ws.Cells.Item(99, 11).Value = 1000
$ws.Cells.Item(99, 13).Value = 498
$ws.Cells.Item(107, 8).Value = 1000
$ws.Cells.Item(107, 9).Value = 1000
$ws.Cells.Item(107, 11).Value = 1000
$ws.Cells.Item(107, 13).Value = 920
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5680.2144
$ws.Cells.Item(31, 9).Value = 2554.6667
$ws.Cells.Item(31, 10).Value = 6532.636
$ws.Cells.Item(31, 11).Value = 2554.6667
$ws.Cells.Item(31, 12).Value = 6532.636
$ws.Cells.Item(31, 13).Value = -2259.6667
$ws.Cells.Item(31, 14).Value = -7122.636
$ws.Cells.Item(34, 8).Value = 5680.2144
$ws.Cells.Item(34, 9).Value = 2554.6667
$ws.Cells.Item(34, 10).Value = 6532.636
$ws.Cells.Item(34, 11).Value = 2554.6667
$ws.Cells.Item(34, 12).Value = 6532.636
$ws.Cells.Item(34, 13).Value = -2352.6667
$ws.Cells.Item(34, 14).Value = -6936.636
$ws.Cells.Item(68, 8).Value = 31664.666
$ws.Cells.Item(68, 9).Value = 24995
$ws.Cells.Item(68, 10).Value = 34999.5
$ws.Cells.Item(68, 11).Value = 24995
$ws.Cells.Item(68, 12).Value = 34999.5
$ws.Cells.Item(68, 13).Value = -24246
$ws.Cells.Item(68, 14).Value = -36497.5
$ws.Cells.Item(71, 8).Value = 31664.666
$ws.Cells.Item(71, 9).Value = 24995
$ws.Cells.Item(71, 10).Value = 34999.5
$ws.Cells.Item(71, 11).Value = 74985
$ws.Cells.Item(71, 12).Value = 104998.5
$ws.Cells.Item(71, 13).Value = -71241
$ws.Cells.Item(71, 14).Value = -112486.5
$ws.Cells.Item(105, 8).Value = 2930.4
$ws.Cells.Item(105, 10).Value = 2984.25
$ws.Cells.Item(105, 12).Value = 2984.25
$ws.Cells.Item(105, 14).Value = -6478.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 3321767.8
$ws.Cells.Item(32, 9).Value = 2713.1428
$ws.Cells.Item(32, 11).Value = 8139.428400000001
$ws.Cells.Item(32, 13).Value = -7856.428400000001
$ws.Cells.Item(75, 8).Value = 410.33334
$ws.Cells.Item(75, 9).Value = 410.33334
$ws.Cells.Item(75, 11).Value = 1231.00002
$ws.Cells.Item(75, 13).Value = -233.0000199999999
$ws.Cells.Item(78, 8).Value = 410.33334
$ws.Cells.Item(78, 9).Value = 410.33334
$ws.Cells.Item(78, 11).Value = 3693.00006
$ws.Cells.Item(78, 13).Value = 1298.99994
$ws.Cells.Item(117, 8).Value = 1852.2727
$ws.Cells.Item(117, 9).Value = 933.3333
$ws.Cells.Item(117, 10).Value = 2196.875
$ws.Cells.Item(117, 11).Value = 2799.9999
$ws.Cells.Item(117, 12).Value = 6590.625
$ws.Cells.Item(117, 13).Value = 642.0001000000002
$ws.Cells.Item(117, 14).Value = -13474.625
$ws.Cells.Item(121, 8).Value = 764.44446
$ws.Cells.Item(121, 9).Value = 206
$ws.Cells.Item(121, 10).Value = 1462.5
$ws.Cells.Item(121, 11).Value = 618
$ws.Cells.Item(121, 12).Value = 4387.5
$ws.Cells.Item(121, 13).Value = 692
$ws.Cells.Item(121, 14).Value = -7007.5
$ws.Cells.Item(122, 8).Value = 929.2222
$ws.Cells.Item(122, 10).Value = 1058.2
$ws.Cells.Item(122, 12).Value = 9523.800000000001
$ws.Cells.Item(122, 14).Value = -14423.8
$ws.Cells.Item(129, 8).Value = 7944.9
$ws.Cells.Item(129, 9).Value = 2149.6667
$ws.Cells.Item(129, 10).Value = 10428.571
$ws.Cells.Item(129, 11).Value = 6449.000100000001
$ws.Cells.Item(129, 12).Value = 31285.713
$ws.Cells.Item(129, 13).Value = -1449.000100000001
$ws.Cells.Item(129, 14).Value = -41285.713
$ws.Cells.Item(131, 8).Value = 4121.913
$ws.Cells.Item(131, 9).Value = 2722.818
$ws.Cells.Item(131, 10).Value = 5404.4165
$ws.Cells.Item(131, 11).Value = 8168.454000000001
$ws.Cells.Item(131, 12).Value = 16213.2495
$ws.Cells.Item(131, 13).Value = -3128.454000000001
$ws.Cells.Item(131, 14).Value = -26293.2495
$ws.Cells.Item(134, 8).Value = 3723
$ws.Cells.Item(134, 10).Value = 13999
$ws.Cells.Item(134, 12).Value = 41997
$ws.Cells.Item(134, 14).Value = -52137
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1984.88
$ws.Cells.Item(97, 9).Value = 1935.65
$ws.Cells.Item(97, 11).Value = 1935.65
$ws.Cells.Item(97, 13).Value = -1439.65
$ws.Cells.Item(126, 8).Value = 3730.6155
$ws.Cells.Item(126, 9).Value = 1899.5
$ws.Cells.Item(126, 10).Value = 4544.4443
$ws.Cells.Item(126, 11).Value = 5698.5
$ws.Cells.Item(126, 12).Value = 13633.3329
$ws.Cells.Item(126, 13).Value = -3228.5
$ws.Cells.Item(126, 14).Value = -18573.3329
$ws.Cells.Item(136, 8).Value = 26711.5
$ws.Cells.Item(136, 10).Value = 26711.5
$ws.Cells.Item(136, 12).Value = 80134.5
$ws.Cells.Item(136, 14).Value = -85234.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 418.7647
$ws.Cells.Item(55, 9).Value = 375.30768
$ws.Cells.Item(55, 11).Value = 375.30768
$ws.Cells.Item(55, 13).Value = -202.30768
$ws.Cells.Item(132, 8).Value = 5254.091
$ws.Cells.Item(132, 9).Value = 5254.091
$ws.Cells.Item(132, 11).Value = 15762.273
$ws.Cells.Item(132, 13).Value = -13232.273
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1934.3334
$ws.Cells.Item(122, 9).Value = 1934.3334
$ws.Cells.Item(122, 11).Value = 5803.0002
$ws.Cells.Item(122, 13).Value = -3353.0002
$ws.Cells.Item(132, 8).Value = 1357.3478
$ws.Cells.Item(132, 9).Value = 1391.5454
$ws.Cells.Item(132, 11).Value = 4174.6362
$ws.Cells.Item(132, 13).Value = -1644.6362
$ws.Cells.Item(136, 8).Value = 33597.87
$ws.Cells.Item(136, 9).Value = 1065.65
$ws.Cells.Item(136, 10).Value = 92747.37
$ws.Cells.Item(136, 11).Value = 3196.95
$ws.Cells.Item(136, 12).Value = 278242.11
$ws.Cells.Item(136, 13).Value = -646.9500000000003
$ws.Cells.Item(136, 14).Value = -283342.11
